# Updated cryptos list on Thu Oct 10 20:59:21 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns on the active sheet.
# Some Price values happen to parse as plain numbers (e.g. "559.11"); those
# are written with a leading apostrophe so Excel keeps them as text (as in
# the source data, which mixes thousand-dot formatted strings with plain
# decimals), then the style is reset to "Normal" so no stray number-format
# style gets attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.760.68'
$ws.Range("E2").Value = '  -1.20%  '
$ws.Range("D3").Value = '2.370.19'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = "'559.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").Value = "'137.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.89%  '
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").Value = "'0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("D9").Value = '2.366.46'
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("E10").Value = '  -1.66%  '
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("D14").Value = "'25.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("D15").Value = '2.794.40'
$ws.Range("D16").Value = "'0.0000165"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.01%  '
$ws.Range("D17").Value = '59.694.25'
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("D18").Value = '2.374.53'
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("D19").Value = "'7.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +11.22%  '
$ws.Range("D20").Value = "'10.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("D21").Value = "'321.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("E22").Value = '  +1.76%  '
$ws.Range("D23").Value = "'6.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.56%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  -2.97%  '
$ws.Range("D26").Value = "'64.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").Value = "'559.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("E28").Value = '  -6.30%  '
$ws.Range("D29").Value = '2.483.01'
$ws.Range("E29").Value = '  -1.12%  '
$ws.Range("D30").Value = '0.0₃0922'
$ws.Range("E30").Value = '  +2.10%  '
$ws.Range("D31").Value = "'8.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.99%  '
$ws.Range("E32").Value = '  -1.87%  '
$ws.Range("E33").Value = '  -2.46%  '
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("E36").Value = '  +4.18%  '
$ws.Range("D37").Value = "'152.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.67%  '
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("D39").Value = "'4.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("D40").Value = "'18.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("D41").Value = "'4.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.75%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").Value = "'41.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("E45").Value = '  +3.72%  '
$ws.Range("D46").Value = '0.0₆0297'
$ws.Range("E46").Value = '  +6.69%  '
$ws.Range("D47").Value = "'138.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("E48").Value = '  +1.19%  '
$ws.Range("D49").Value = "'0.585"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("E51").Value = '  -0.69%  '
